$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 418.75
$ws.Range("I6").Value = 50
$ws.Range("J6").Value = 471.42856
$ws.Range("K6").Value = 150
$ws.Range("L6").Value = 1414.28568
$ws.Range("M6").Value = -38
$ws.Range("N6").Value = -1638.28568
$ws.Range("H8").Value = 36.333332
$ws.Range("I8").Value = 36.333332
$ws.Range("K8").Value = 108.999996
$ws.Range("M8").Value = 30.000004
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H41").Value = 550
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 550
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 550
$ws.Range("N41").Value = -1430
$ws.Range("M41").ClearContents()
$ws.Range("H42").Value = 83.09999999999999
$ws.Range("I42").Value = 36.333332
$ws.Range("J42").Value = 103.14286
$ws.Range("K42").Value = 108.999996
$ws.Range("L42").Value = 309.42858
$ws.Range("M42").Value = 121.000004
$ws.Range("N42").Value = -769.42858
$ws.Range("H116").Value = 2427.6206
$ws.Range("I116").Value = 2324.1538
$ws.Range("K116").Value = 2324.1538
$ws.Range("M116").Value = 1117.8462
$ws.Range("H129").Value = 986.325
$ws.Range("I129").Value = 914.0769
$ws.Range("J129").Value = 1000.34326
$ws.Range("K129").Value = 2742.2307
$ws.Range("L129").Value = 3001.02978
$ws.Range("M129").Value = 2257.7693
$ws.Range("N129").Value = -13001.02978
$ws.Range("H141").Value = 3074.6453
$ws.Range("I141").Value = 1175.4166
$ws.Range("J141").Value = 9586.286
$ws.Range("K141").Value = 3526.2498
$ws.Range("L141").Value = 28758.858
$ws.Range("M141").Value = 1653.7502
$ws.Range("N141").Value = -39118.858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1143.5834
$ws.Range("I2").Value = 1082.5555
$ws.Range("J2").Value = 1326.6666
$ws.Range("K2").Value = 1082.5555
$ws.Range("L2").Value = 1326.6666
$ws.Range("M2").Value = -969.5554999999999
$ws.Range("N2").Value = -1552.6666
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("H32").Value = 24396696
$ws.Range("I32").Value = 6693.4287
$ws.Range("K32").Value = 6693.4287
$ws.Range("M32").Value = -6406.4287
$ws.Range("H45").Value = 92407.17999999999
$ws.Range("I45").Value = 251045
$ws.Range("J45").Value = 1757
$ws.Range("K45").Value = 251045
$ws.Range("L45").Value = 1757
$ws.Range("M45").Value = -250668
$ws.Range("N45").Value = -2511
$ws.Range("H116").Value = 1143.5834
$ws.Range("I116").Value = 1082.5555
$ws.Range("J116").Value = 1326.6666
$ws.Range("K116").Value = 1082.5555
$ws.Range("L116").Value = 1326.6666
$ws.Range("M116").Value = 1211.4445
$ws.Range("N116").Value = -5914.6666
$ws.Range("H132").Value = 920109.2
$ws.Range("I132").Value = 898.5741
$ws.Range("J132").Value = 5883846.5
$ws.Range("K132").Value = 2695.7223
$ws.Range("L132").Value = 17651539.5
$ws.Range("M132").Value = -165.7223000000004
$ws.Range("N132").Value = -17656599.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1143.5834
$ws.Range("I3").Value = 1082.5555
$ws.Range("J3").Value = 1326.6666
$ws.Range("K3").Value = 1082.5555
$ws.Range("L3").Value = 1326.6666
$ws.Range("M3").Value = -968.5554999999999
$ws.Range("N3").Value = -1554.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 591.3
$ws.Range("I107").Value = 445.3684
$ws.Range("J107").Value = 843.36365
$ws.Range("K107").Value = 445.3684
$ws.Range("L107").Value = 843.36365
$ws.Range("M107").Value = 1474.6316
$ws.Range("N107").Value = -4683.36365

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 157172500
$ws.Range("J9").Value = 166701250
$ws.Range("L9").Value = 500103750
$ws.Range("N9").Value = -500104198
$ws.Range("H10").Value = 145.33333
$ws.Range("I10").Value = 114.4
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 343.2
$ws.Range("L10").Value = 900
$ws.Range("M10").Value = -204.2
$ws.Range("N10").Value = -1178
$ws.Range("H17").Value = 367.66666
$ws.Range("I17").Value = 76.5
$ws.Range("K17").Value = 229.5
$ws.Range("M17").Value = -60.5
$ws.Range("H19").Value = 1800
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1800
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 5400
$ws.Range("N19").Value = -5748
$ws.Range("M19").ClearContents()
$ws.Range("H46").Value = 914.2857
$ws.Range("I46").Value = 350
$ws.Range("J46").Value = 1666.6666
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 4999.9998
$ws.Range("M46").Value = -959
$ws.Range("N46").Value = -5181.9998
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15500
$ws.Range("H92").Value = 6773.722
$ws.Range("I92").Value = 764.6667
$ws.Range("J92").Value = 7975.533
$ws.Range("K92").Value = 2294.0001
$ws.Range("L92").Value = 23926.599
$ws.Range("M92").Value = -1046.0001
$ws.Range("N92").Value = -26422.599
$ws.Range("H131").Value = 817.9091
$ws.Range("J131").Value = 824.47424
$ws.Range("L131").Value = 2473.42272
$ws.Range("N131").Value = -12553.42272

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2600
$ws.Range("J9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("N9").Value = -5340
$ws.Range("H17").Value = 3000
$ws.Range("J17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("N17").Value = -3336
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H64").Value = 19800
$ws.Range("J64").Value = 19800
$ws.Range("L64").Value = 19800
$ws.Range("N64").Value = -20296
$ws.Range("H67").Value = 19800
$ws.Range("J67").Value = 19800
$ws.Range("L67").Value = 19800
$ws.Range("N67").Value = -21516
$ws.Range("H109").Value = 20162
$ws.Range("J109").Value = 20162
$ws.Range("L109").Value = 20162
$ws.Range("N109").Value = -22242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 499
$ws.Range("I12").Value = 498
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 498
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -328
$ws.Range("N12").Value = -840
$ws.Range("H16").Value = 107143250
$ws.Range("I16").Value = 10204496
$ws.Range("J16").Value = 333333660
$ws.Range("K16").Value = 10204496
$ws.Range("L16").Value = 333333660
$ws.Range("M16").Value = -10204326
$ws.Range("N16").Value = -333334000
$ws.Range("H19").Value = 250
$ws.Range("I19").Value = 250
$ws.Range("K19").Value = 250
$ws.Range("M19").Value = -80
$ws.Range("H21").Value = 4800
$ws.Range("J21").Value = 4800
$ws.Range("L21").Value = 4800
$ws.Range("N21").Value = -5148
$ws.Range("H25").Value = 1057
$ws.Range("I25").Value = 1057
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1057
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -827
$ws.Range("N25").ClearContents()
$ws.Range("H61").Value = 1305.1333
$ws.Range("I61").Value = 1382.1666
$ws.Range("J61").Value = 1253.7778
$ws.Range("K61").Value = 1382.1666
$ws.Range("L61").Value = 1253.7778
$ws.Range("M61").Value = -1180.1666
$ws.Range("N61").Value = -1657.7778
$ws.Range("H113").Value = 1305.1333
$ws.Range("I113").Value = 1382.1666
$ws.Range("J113").Value = 1253.7778
$ws.Range("K113").Value = 1382.1666
$ws.Range("L113").Value = 1253.7778
$ws.Range("M113").Value = 787.8334
$ws.Range("N113").Value = -5593.7778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 8000
$ws.Range("J18").Value = 8000
$ws.Range("L18").Value = 8000
$ws.Range("N18").Value = -8346
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H132").Value = 19036.54
$ws.Range("I132").Value = 22587.797
$ws.Range("J132").Value = 8160.8125
$ws.Range("K132").Value = 67763.391
$ws.Range("L132").Value = 24482.4375
$ws.Range("M132").Value = -65233.391
$ws.Range("N132").Value = -29542.4375
$ws.Range("H133").Value = 31500
$ws.Range("J133").Value = 31500
$ws.Range("L133").Value = 31500
$ws.Range("N133").Value = -41620
